# Append new activity-log entries (4/6/2019 - 4/9/2019) to the "2019" sheet's
# Table2, matching the source author's manual data-entry session.
#
# Strategy: grow the table by the required number of rows via ListRows.Add()
# (this keeps the ListObject/table ref + worksheet dimension in sync), then
# for each new row copy the cell *formatting* down from the nearest existing
# row of the same "shape" (a Sleep row has Start+End+Activity but no Comment;
# a Food row has Start+Activity+Comment but no End) before writing the real
# values. This mirrors how the values were actually typed into Excel, where
# each new row inherits the formatting of the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$tbl = $ws.ListObjects.Item("Table2")

$DQ = [char]34
$zFormula = "=IF(Table2[[#This Row],[Activity]]=" + $DQ + "Sleep" + $DQ + ",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24," + $DQ + "NA" + $DQ + ")"

# Template rows already present in the sheet:
#   180 -> Start, End, Activity="Sleep"               (no Comment)
#   181 -> Start, Activity="Food", Comment             (no End)
$sleepTemplateRow = 180
$foodTemplateRow  = 181

$newRows = @(
    @{ Row = 191; Start = 43561.941469907404; End = 43562.167361111111; Activity = "Sleep"; Comment = $null },
    @{ Row = 192; Start = 43561.8125;          End = $null;             Activity = "Food";  Comment = "Pizza + breaded chicken" },
    @{ Row = 193; Start = 43562.180555555555;  End = $null;             Activity = "Food";  Comment = "Latte" },
    @{ Row = 194; Start = 43561.520833333336;  End = $null;             Activity = "Food";  Comment = "Whole wheat pasta" },
    @{ Row = 195; Start = 43562.40625;         End = $null;             Activity = "Food";  Comment = "eggs + cheese" },
    @{ Row = 196; Start = 43563.40625;         End = $null;             Activity = "Food";  Comment = "egg/banana/fiber" },
    @{ Row = 197; Start = 43563.791666666664;  End = $null;             Activity = "Food";  Comment = "Beans + cabbage" },
    @{ Row = 198; Start = 43563.583333333336;  End = $null;             Activity = "Food";  Comment = "Clam Chowder" },
    @{ Row = 199; Start = 43562.939583333333;  End = 43563.228472222225; Activity = "Sleep"; Comment = $null },
    @{ Row = 200; Start = 43563.90761574074;   End = 43564.199305555558; Activity = "Sleep"; Comment = $null },
    @{ Row = 201; Start = 43564.21875;         End = $null;             Activity = "Food";  Comment = "Latte" }
)

foreach ($entry in $newRows) {
    # Grow the table by one row (keeps table ref / autofilter / dimension in sync).
    $tbl.ListRows.Add() | Out-Null

    $r = $entry.Row
    $isSleep = ($entry.Activity -eq "Sleep")
    $templateRow = $(if ($isSleep) { $sleepTemplateRow } else { $foodTemplateRow })

    # Copy the whole row's formatting down onto the freshly added row.
    $ws.Range("A" + $templateRow + ":E" + $templateRow).Copy()
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $entry.Start

    if ($entry.End -ne $null) {
        $ws.Cells.Item($r, 2).Value = $entry.End
    } else {
        $ws.Cells.Item($r, 2).Clear()
    }

    $ws.Cells.Item($r, 3).Value = $entry.Activity

    if ($entry.Comment -ne $null) {
        $ws.Cells.Item($r, 4).Value = $entry.Comment
    } else {
        $ws.Cells.Item($r, 4).Clear()
    }

    $ws.Cells.Item($r, 5).Formula = $zFormula
}

$excel.Application.CutCopyMode = $false
